# Round 6 update: append 3 new match columns (JW, JX, JY) to the stats
# table. They carry the round-6 figures for each stat row; all three
# columns end up identical (matches the source data's duplication).
# Column JV (previously the last/rightmost column) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# JW=283, JX=284, JY=285 (1-indexed columns)
$newVals = @{
    1 = 10800
    2 = 2023
    3 = 6
    4 = 1
    5 = 1
    6 = 79
    7 = 76
    8 = 3
    9 = 1
    10 = 10
    11 = 205
    12 = 150
    13 = 355
    14 = 1.37
    15 = 79
    16 = 75
    17 = 45
    18 = 20
    19 = 19
    20 = 11
    21 = 8
    22 = 8
    23 = 5
    24 = 24
    25 = 45.8
    26 = 32.27
    27 = 14.79
    28 = 47
    29 = 63
    30 = 34
    31 = 53
    32 = 2.21
    33 = 4.82
    34 = 35.8
    35 = 20.8
    36 = 187.3
    37 = 87.59999999999999
    38 = 24.74
    39 = 82.59999999999999
    40 = 9
    41 = 10
    42 = 0
    43 = 4
    44 = 145
    45 = 206
    46 = 255
    47 = 71.8
    48 = 63
    49 = 7
    50 = 14
    51 = 47
    52 = 34
    53 = 43
    54 = 3
    55 = 8
    56 = 72.7
    57 = 185
    58 = 167
    59 = 352
    60 = 1.11
    61 = 78
    62 = 88
    63 = 52
    64 = 19
    65 = 20
    66 = 11
    67 = 9
    68 = 9
    69 = 1
    70 = 21
    71 = 52.4
    72 = 32
    73 = 16.76
    74 = 38
    75 = 54
    76 = 42
    77 = 47
    78 = 2.24
    79 = 4.27
    80 = 42.6
    81 = 23.4
    82 = 187.7
    83 = 85.3
    84 = 24.58
    85 = 74.2
    86 = 11
    87 = 5
    88 = 4
    89 = 3
    90 = 134
    91 = 212
    92 = 265
    93 = 75.3
    94 = 54
    95 = 8
    96 = 7
    97 = 38
    98 = 42
    99 = 42
    100 = 6
    101 = 9
    102 = 81.8
}

foreach ($r in 1..102) {
    $v = $newVals[$r]
    $ws.Cells.Item($r, 283).Value = $v   # JW<r>
    $ws.Cells.Item($r, 284).Value = $v   # JX<r>
    $ws.Cells.Item($r, 285).Value = $v   # JY<r>
}